# Daily attendance processing - 2026-01-28 11:12:50
# For every "Recorded By" (column G) cell whose value is a comma-separated
# list of recorders, rotate the list left by one position (move the first
# entry to the end), except for the literal value "admin@admin.com, System"
# which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$skipValue = "admin@admin.com, System"
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    if ($text -eq $skipValue) {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -le 1) {
        continue
    }

    $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
    $newText = $rotated -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
        $changed++
    }
}

Write-Output "Rotated $changed 'Recorded By' cells"
